$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @('D2', '60.926.79'),
    @('E2', '  -1.94%  '),
    @('D3', '3.386.13'),
    @('E3', '  -0.77%  '),
    @('E4', '  +0.01%  '),
    @('D5', '572.85'),
    @('E5', '  -0.95%  '),
    @('D6', '137.82'),
    @('E6', '  +0.08%  '),
    @('E7', '  +0.05%  '),
    @('D8', '3.386.26'),
    @('E8', '  -0.73%  '),
    @('E9', '  -1.29%  '),
    @('E10', '  +2.00%  '),
    @('D11', '0.123'),
    @('E11', '  -2.49%  '),
    @('D12', '0.384'),
    @('E12', '  -2.34%  '),
    @('D13', '3.965.08'),
    @('E13', '  -0.66%  '),
    @('D15', '26.01'),
    @('E15', '  +1.80%  '),
    @('E16', '  -3.17%  '),
    @('D17', '3.395.56'),
    @('E17', '  -0.55%  '),
    @('D18', '61.083.07'),
    @('E18', '  -1.73%  '),
    @('D19', '13.98'),
    @('E19', '  -1.85%  '),
    @('D20', '5.80'),
    @('E20', '  -1.34%  '),
    @('D21', '9.43'),
    @('E21', '  -0.74%  '),
    @('D22', '375.74'),
    @('E22', '  -3.58%  '),
    @('D23', '0.554'),
    @('E23', '  -2.99%  '),
    @('D24', '3.529.47'),
    @('E24', '  -0.46%  '),
    @('E25', '  -0.17%  '),
    @('D26', '0.0000125'),
    @('E26', '  -2.55%  '),
    @('D27', '71.09'),
    @('E27', '  -0.62%  '),
    @('D28', '1.76'),
    @('E28', '  +10.89%  '),
    @('E29', '  +9.17%  '),
    @('D30', '7.51'),
    @('E30', '  -2.95%  '),
    @('E31', '  +0.02%  '),
    @('D32', '8.15'),
    @('E32', '  -1.84%  '),
    @('D33', '2.15'),
    @('E33', '  -1.29%  '),
    @('E34', '  +0.00%  '),
    @('E35', '  +0.42%  '),
    @('D36', '5.21'),
    @('E36', '  -4.34%  '),
    @('D37', '6.87'),
    @('E37', '  -1.94%  '),
    @('E38', '  -0.86%  '),
    @('D39', '164.74'),
    @('E39', '  +0.70%  '),
    @('E41', '  +0.01%  '),
    @('D42', '0.779'),
    @('E42', '  -0.97%  '),
    @('E43', '  -4.90%  '),
    @('D44', '41.65'),
    @('E44', '  -0.13%  '),
    @('E45', '  -1.51%  '),
    @('E46', '  -3.28%  '),
    @('D47', '24.32'),
    @('E47', '  -3.42%  '),
    @('D48', '2.453.44'),
    @('E48', '  +3.32%  '),
    @('D49', '23.10'),
    @('E49', '  -1.37%  '),
    @('E50', '  -2.53%  '),
    @('D51', '2.41'),
    @('E51', '  +4.25%  ')
)

foreach ($pair in $pairs) {
    $cellRef = $pair[0]
    $val = $pair[1]
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}
